# "Scrutinizing Carnagey and Anderson 2005 entry."
#
# The sheet already has an AutoFilter (A1:AJ85) restricting column X
# (setting2) to Exp/Nonexp. We narrow the view further by also filtering
# column O (Full.Reference) down to just the Carnagey & Anderson (2005)
# study, then re-sort the now-small set of matching rows by column R
# (Study), and move the selection onto the newly focused data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$carnagey = "Carnagey, N. L., & Anderson, C.A. (2005). The effects of reward and punishment in violent video games on aggressive affect, cognition, and behavior. Psychological Science, 16, 882-889."

$fullRange = $ws.Range("A1:AJ85")

# Field 15 = column O (Full.Reference): keep only the Carnagey & Anderson row.
$fullRange.AutoFilter(15, @($carnagey), 7)

# Re-apply the pre-existing field 24 = column X (setting2) filter so both
# filter columns are (re)written together, in field order.
$fullRange.AutoFilter(24, @("Exp", "Nonexp"), 7)

# The visible rows (34-42) all belong to the Carnagey & Anderson study;
# sort just that block by column R ("Study"), ascending.
$ws.Range("A34:AJ42").Sort($ws.Range("R34:R42"), 1)

# Reflect the new focus: selection on J40, scrolled back to the top.
$ws.Range("J40").Select()
